# Rename the existing sheet to "NYT Books", add a new "Box Office" sheet
# right after it, populate it with the scraped daily box-office data, and
# format it to match the rest of the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet1 -> "NYT Books" -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "NYT Books"

# --- New "Box Office" sheet, placed right after "NYT Books" ----------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Box Office"

# --- Header row --------------------------------------------------------------
$ws2.Range("A1").Value = "Date"
$ws2.Range("B1").Value = "Gross"

# --- Daily date / gross data, rows 2-184 --------------------------------------
$data = New-Object 'object[,]' 183,2
$data[0,0] = 41998
$data[0,1] = 240211
$data[1,0] = 41999
$data[1,1] = 199481
$data[2,0] = 42000
$data[2,1] = 215858
$data[3,0] = 42001
$data[3,1] = 218117
$data[4,0] = 42002
$data[4,1] = 166444
$data[5,0] = 42003
$data[5,1] = 164201
$data[6,0] = 42004
$data[6,1] = 124987
$data[7,0] = 42005
$data[7,1] = 223080
$data[8,0] = 42006
$data[8,1] = 233719
$data[9,0] = 42007
$data[9,1] = 238862
$data[10,0] = 42008
$data[10,1] = 204328
$data[11,0] = 42009
$data[11,1] = 100505
$data[12,0] = 42010
$data[12,1] = 89058
$data[13,0] = 42011
$data[13,1] = 93332
$data[14,0] = 42012
$data[14,1] = 82968
$data[15,0] = 42013
$data[15,1] = 156360
$data[16,0] = 42014
$data[16,1] = 234210
$data[17,0] = 42015
$data[17,1] = 188948
$data[18,0] = 42016
$data[18,1] = 60183
$data[19,0] = 42017
$data[19,1] = 68600
$data[20,0] = 42018
$data[20,1] = 69270
$data[21,0] = 42019
$data[21,1] = 52056
$data[22,0] = 42020
$data[22,1] = 30338488
$data[23,0] = 42021
$data[23,1] = 34547284
$data[24,0] = 42022
$data[24,1] = 24383294
$data[25,0] = 42023
$data[25,1] = 17942391
$data[26,0] = 42024
$data[26,1] = 9924117
$data[27,0] = 42025
$data[27,1] = 7555269
$data[28,0] = 42026
$data[28,1] = 7656492
$data[29,0] = 42027
$data[29,1] = 18213554
$data[30,0] = 42028
$data[30,1] = 28635135
$data[31,0] = 42029
$data[31,1] = 17779615
$data[32,0] = 42030
$data[32,1] = 4183367
$data[33,0] = 42031
$data[33,1] = 5017038
$data[34,0] = 42032
$data[34,1] = 3823135
$data[35,0] = 42033
$data[35,1] = 3668056
$data[36,0] = 42034
$data[36,1] = 9905616
$data[37,0] = 42035
$data[37,1] = 16510536
$data[38,0] = 42036
$data[38,1] = 4244376
$data[39,0] = 42037
$data[39,1] = 2645109
$data[40,0] = 42038
$data[40,1] = 2923141
$data[41,0] = 42039
$data[41,1] = 2273342
$data[42,0] = 42040
$data[42,1] = 2506106
$data[43,0] = 42041
$data[43,1] = 6163365
$data[44,0] = 42042
$data[44,1] = 11032447
$data[45,0] = 42043
$data[45,1] = 6093301
$data[46,0] = 42044
$data[46,1] = 1590242
$data[47,0] = 42045
$data[47,1] = 1773361
$data[48,0] = 42046
$data[48,1] = 1468160
$data[49,0] = 42047
$data[49,1] = 1477178
$data[50,0] = 42048
$data[50,1] = 3745563
$data[51,0] = 42049
$data[51,1] = 7824072
$data[52,0] = 42050
$data[52,1] = 4845170
$data[53,0] = 42051
$data[53,1] = 2365038
$data[54,0] = 42052
$data[54,1] = 1284385
$data[55,0] = 42053
$data[55,1] = 1037497
$data[56,0] = 42054
$data[56,1] = 1157160
$data[57,0] = 42055
$data[57,1] = 2718122
$data[58,0] = 42056
$data[58,1] = 4595743
$data[59,0] = 42057
$data[59,1] = 2738582
$data[60,0] = 42058
$data[60,1] = 836221
$data[61,0] = 42059
$data[61,1] = 948496
$data[62,0] = 42060
$data[62,1] = 792101
$data[63,0] = 42061
$data[63,1] = 822004
$data[64,0] = 42062
$data[64,1] = 1903469
$data[65,0] = 42063
$data[65,1] = 3601727
$data[66,0] = 42064
$data[66,1] = 1889097
$data[67,0] = 42065
$data[67,1] = 502179
$data[68,0] = 42066
$data[68,1] = 507318
$data[69,0] = 42067
$data[69,1] = 429402
$data[70,0] = 42068
$data[70,1] = 457340
$data[71,0] = 42069
$data[71,1] = 1146840
$data[72,0] = 42070
$data[72,1] = 2158539
$data[73,0] = 42071
$data[73,1] = 1096012
$data[74,0] = 42072
$data[74,1] = 342258
$data[75,0] = 42073
$data[75,1] = 409581
$data[76,0] = 42074
$data[76,1] = 363324
$data[77,0] = 42075
$data[77,1] = 354030
$data[78,0] = 42076
$data[78,1] = 755768
$data[79,0] = 42077
$data[79,1] = 1382406
$data[80,0] = 42078
$data[80,1] = 673168
$data[81,0] = 42079
$data[81,1] = 232068
$data[82,0] = 42080
$data[82,1] = 255383
$data[83,0] = 42081
$data[83,1] = 256310
$data[84,0] = 42082
$data[84,1] = 236086
$data[85,0] = 42083
$data[85,1] = 432738
$data[86,0] = 42084
$data[86,1] = 785413
$data[87,0] = 42085
$data[87,1] = 390100
$data[88,0] = 42086
$data[88,1] = 152393
$data[89,0] = 42087
$data[89,1] = 155245
$data[90,0] = 42088
$data[90,1] = 137446
$data[91,0] = 42089
$data[91,1] = 126147
$data[92,0] = 42090
$data[92,1] = 246937
$data[93,0] = 42091
$data[93,1] = 423874
$data[94,0] = 42092
$data[94,1] = 203242
$data[95,0] = 42093
$data[95,1] = 81213
$data[96,0] = 42094
$data[96,1] = 91324
$data[97,0] = 42095
$data[97,1] = 83249
$data[98,0] = 42096
$data[98,1] = 88355
$data[99,0] = 42097
$data[99,1] = 230878
$data[100,0] = 42098
$data[100,1] = 261490
$data[101,0] = 42099
$data[101,1] = 135692
$data[102,0] = 42100
$data[102,1] = 55473
$data[103,0] = 42101
$data[103,1] = 61575
$data[104,0] = 42102
$data[104,1] = 55517
$data[105,0] = 42103
$data[105,1] = 57183
$data[106,0] = 42104
$data[106,1] = 180278
$data[107,0] = 42105
$data[107,1] = 302268
$data[108,0] = 42106
$data[108,1] = 112495
$data[109,0] = 42107
$data[109,1] = 40555
$data[110,0] = 42108
$data[110,1] = 48387
$data[111,0] = 42109
$data[111,1] = 44459
$data[112,0] = 42110
$data[112,1] = 38505
$data[113,0] = 42111
$data[113,1] = 104872
$data[114,0] = 42112
$data[114,1] = 185191
$data[115,0] = 42113
$data[115,1] = 89655
$data[116,0] = 42114
$data[116,1] = 37173
$data[117,0] = 42115
$data[117,1] = 42393
$data[118,0] = 42116
$data[118,1] = 37243
$data[119,0] = 42117
$data[119,1] = 40967
$data[120,0] = 42118
$data[120,1] = 96576
$data[121,0] = 42119
$data[121,1] = 187408
$data[122,0] = 42120
$data[122,1] = 90615
$data[123,0] = 42121
$data[123,1] = 31715
$data[124,0] = 42122
$data[124,1] = 34857
$data[125,0] = 42123
$data[125,1] = 29628
$data[126,0] = 42124
$data[126,1] = 26072
$data[127,0] = 42125
$data[127,1] = 123899
$data[128,0] = 42126
$data[128,1] = 171757
$data[129,0] = 42127
$data[129,1] = 148259
$data[130,0] = 42128
$data[130,1] = 45960
$data[131,0] = 42129
$data[131,1] = 50556
$data[132,0] = 42130
$data[132,1] = 39939
$data[133,0] = 42131
$data[133,1] = 35947
$data[134,0] = 42132
$data[134,1] = 86330
$data[135,0] = 42133
$data[135,1] = 119155
$data[136,0] = 42134
$data[136,1] = 85256
$data[137,0] = 42135
$data[137,1] = 34102
$data[138,0] = 42136
$data[138,1] = 37513
$data[139,0] = 42137
$data[139,1] = 28134
$data[140,0] = 42138
$data[140,1] = 27572
$data[141,0] = 42139
$data[141,1] = 54584
$data[142,0] = 42140
$data[142,1] = 88711
$data[143,0] = 42141
$data[143,1] = 45686
$data[144,0] = 42142
$data[144,1] = 15890
$data[145,0] = 42143
$data[145,1] = 17002
$data[146,0] = 42144
$data[146,1] = 16151
$data[147,0] = 42145
$data[147,1] = 19120
$data[148,0] = 42146
$data[148,1] = 43350
$data[149,0] = 42147
$data[149,1] = 69227
$data[150,0] = 42148
$data[150,1] = 72440
$data[151,0] = 42149
$data[151,1] = 42529
$data[152,0] = 42150
$data[152,1] = 13184
$data[153,0] = 42151
$data[153,1] = 10942
$data[154,0] = 42152
$data[154,1] = 12364
$data[155,0] = 42153
$data[155,1] = 19029
$data[156,0] = 42154
$data[156,1] = 28047
$data[157,0] = 42155
$data[157,1] = 20327
$data[158,0] = 42156
$data[158,1] = 6327
$data[159,0] = 42157
$data[159,1] = 7076
$data[160,0] = 42158
$data[160,1] = 6651
$data[161,0] = 42159
$data[161,1] = 6317
$data[162,0] = 42160
$data[162,1] = 6531
$data[163,0] = 42161
$data[163,1] = 10859
$data[164,0] = 42162
$data[164,1] = 5757
$data[165,0] = 42163
$data[165,1] = 2497
$data[166,0] = 42164
$data[166,1] = 2521
$data[167,0] = 42165
$data[167,1] = 2651
$data[168,0] = 42166
$data[168,1] = 2846
$data[169,0] = 42167
$data[169,1] = 3128
$data[170,0] = 42168
$data[170,1] = 5228
$data[171,0] = 42169
$data[171,1] = 3292
$data[172,0] = 42170
$data[172,1] = 1362
$data[173,0] = 42171
$data[173,1] = 1536
$data[174,0] = 42172
$data[174,1] = 1902
$data[175,0] = 42173
$data[175,1] = 1910
$data[176,0] = 42174
$data[176,1] = 1123
$data[177,0] = 42175
$data[177,1] = 1814
$data[178,0] = 42176
$data[178,1] = 1714
$data[179,0] = 42177
$data[179,1] = 607
$data[180,0] = 42178
$data[180,1] = 789
$data[181,0] = 42179
$data[181,1] = 707
$data[182,0] = 42180
$data[182,1] = 789
$ws2.Range("A2:B184").Value = $data

# --- Number formats ------------------------------------------------------------
# Column A: same custom date format used on the "NYT Books" sheet
$ws2.Columns.Item(1).NumberFormat = "yyyy\-mm\-dd;@"
# Column B: currency format ($#,##0)
$ws2.Range("B2:B184").NumberFormat = '"$"#,##0_);[Red]\("$"#,##0\)'

# --- Column widths (best fit) ---------------------------------------------------
$ws2.Columns.Item(2).ColumnWidth = 11

# --- Selection on the new active sheet -------------------------------------------
[void]$ws2.Range("B2").Select()

Write-Host "done"
